# Add data for 2025-07-11
# Updates the 2025 year-to-date (column L) violent crime counts
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = 'Citywide Totals'; Cell = 'L2'; Value = 3472 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L3'; Value = 3618 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L4'; Value = 903 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L5'; Value = 216 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L6'; Value = 3169 }
    @{ Sheet = 'Citywide Totals'; Cell = 'L7'; Value = 11378 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L4'; Value = 44 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L5'; Value = 43 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L6'; Value = 88 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L7'; Value = 383 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L8'; Value = 739 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L19'; Value = 319 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L20'; Value = 290 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L29'; Value = 613 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L31'; Value = 111 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L33'; Value = 536 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L36'; Value = 154 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L37'; Value = 405 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L40'; Value = 31 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L46'; Value = 25 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L47'; Value = 82 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L51'; Value = 142 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L54'; Value = 240 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L55'; Value = 108 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L60'; Value = 67 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L63'; Value = 38 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L65'; Value = 220 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L67'; Value = 405 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L79'; Value = 295 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L81'; Value = 12 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L84'; Value = 113 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L85'; Value = 576 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L91'; Value = 161 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L94'; Value = 136 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L95'; Value = 159 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L97'; Value = 96 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L99'; Value = 190 }
    @{ Sheet = 'By Neighborhood'; Cell = 'L101'; Value = 11378 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'L2'; Value = 126 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'L3'; Value = 115 }
    @{ Sheet = 'Auburn Gresham'; Cell = 'L7'; Value = 383 }
    @{ Sheet = 'South Shore'; Cell = 'L2'; Value = 166 }
    @{ Sheet = 'South Shore'; Cell = 'L3'; Value = 237 }
    @{ Sheet = 'South Shore'; Cell = 'L7'; Value = 576 }
    @{ Sheet = 'Austin'; Cell = 'L2'; Value = 215 }
    @{ Sheet = 'Austin'; Cell = 'L3'; Value = 244 }
    @{ Sheet = 'Austin'; Cell = 'L7'; Value = 739 }
    @{ Sheet = 'Garfield Park'; Cell = 'L2'; Value = 150 }
    @{ Sheet = 'Garfield Park'; Cell = 'L3'; Value = 168 }
    @{ Sheet = 'Garfield Park'; Cell = 'L6'; Value = 180 }
    @{ Sheet = 'Garfield Park'; Cell = 'L7'; Value = 536 }
    @{ Sheet = 'West Pullman'; Cell = 'L3'; Value = 50 }
    @{ Sheet = 'West Pullman'; Cell = 'L7'; Value = 159 }
    @{ Sheet = 'Grand Crossing'; Cell = 'L3'; Value = 127 }
    @{ Sheet = 'Grand Crossing'; Cell = 'L7'; Value = 405 }
    @{ Sheet = 'New City'; Cell = 'L2'; Value = 78 }
    @{ Sheet = 'New City'; Cell = 'L3'; Value = 65 }
    @{ Sheet = 'New City'; Cell = 'L4'; Value = 10 }
    @{ Sheet = 'New City'; Cell = 'L7'; Value = 220 }
    @{ Sheet = 'Woodlawn'; Cell = 'L2'; Value = 51 }
    @{ Sheet = 'Woodlawn'; Cell = 'L4'; Value = 16 }
    @{ Sheet = 'Woodlawn'; Cell = 'L7'; Value = 190 }
    @{ Sheet = 'Gage Park'; Cell = 'L6'; Value = 35 }
    @{ Sheet = 'Gage Park'; Cell = 'L7'; Value = 111 }
    @{ Sheet = 'North Lawndale'; Cell = 'L2'; Value = 118 }
    @{ Sheet = 'North Lawndale'; Cell = 'L3'; Value = 152 }
    @{ Sheet = 'North Lawndale'; Cell = 'L4'; Value = 30 }
    @{ Sheet = 'North Lawndale'; Cell = 'L7'; Value = 405 }
    @{ Sheet = 'South Deering'; Cell = 'L2'; Value = 39 }
    @{ Sheet = 'South Deering'; Cell = 'L7'; Value = 113 }
    @{ Sheet = 'Loop'; Cell = 'L2'; Value = 52 }
    @{ Sheet = 'Loop'; Cell = 'L3'; Value = 56 }
    @{ Sheet = 'Loop'; Cell = 'L6'; Value = 113 }
    @{ Sheet = 'Loop'; Cell = 'L7'; Value = 240 }
    @{ Sheet = 'Englewood'; Cell = 'L2'; Value = 180 }
    @{ Sheet = 'Englewood'; Cell = 'L4'; Value = 31 }
    @{ Sheet = 'Englewood'; Cell = 'L7'; Value = 613 }
    @{ Sheet = 'Chatham'; Cell = 'L3'; Value = 99 }
    @{ Sheet = 'Chatham'; Cell = 'L7'; Value = 319 }
    @{ Sheet = 'Ashburn'; Cell = 'L2'; Value = 37 }
    @{ Sheet = 'Ashburn'; Cell = 'L7'; Value = 88 }
    @{ Sheet = 'Lower West Side'; Cell = 'L3'; Value = 37 }
    @{ Sheet = 'Lower West Side'; Cell = 'L7'; Value = 108 }
    @{ Sheet = 'Jefferson Park'; Cell = 'L2'; Value = 6 }
    @{ Sheet = 'Jefferson Park'; Cell = 'L7'; Value = 25 }
    @{ Sheet = 'Washington Park'; Cell = 'L3'; Value = 67 }
    @{ Sheet = 'Washington Park'; Cell = 'L7'; Value = 161 }
    @{ Sheet = 'Roseland'; Cell = 'L3'; Value = 108 }
    @{ Sheet = 'Roseland'; Cell = 'L7'; Value = 295 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'L3'; Value = 91 }
    @{ Sheet = 'Chicago Lawn'; Cell = 'L7'; Value = 290 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'L2'; Value = 57 }
    @{ Sheet = 'Grand Boulevard'; Cell = 'L7'; Value = 154 }
    @{ Sheet = 'West Loop'; Cell = 'L3'; Value = 31 }
    @{ Sheet = 'West Loop'; Cell = 'L7'; Value = 136 }
    @{ Sheet = 'Kenwood'; Cell = 'L2'; Value = 29 }
    @{ Sheet = 'Kenwood'; Cell = 'L3'; Value = 29 }
    @{ Sheet = 'Kenwood'; Cell = 'L7'; Value = 82 }
    @{ Sheet = 'West Town'; Cell = 'L3'; Value = 19 }
    @{ Sheet = 'West Town'; Cell = 'L7'; Value = 96 }
    @{ Sheet = 'Armour Square'; Cell = 'L5'; Value = 1 }
    @{ Sheet = 'Armour Square'; Cell = 'L7'; Value = 43 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'L2'; Value = 44 }
    @{ Sheet = 'Little Italy, UIC'; Cell = 'L7'; Value = 142 }
    @{ Sheet = 'Morgan Park'; Cell = 'L3'; Value = 27 }
    @{ Sheet = 'Morgan Park'; Cell = 'L7'; Value = 67 }
    @{ Sheet = 'Hegewisch'; Cell = 'L3'; Value = 14 }
    @{ Sheet = 'Hegewisch'; Cell = 'L6'; Value = 7 }
    @{ Sheet = 'Hegewisch'; Cell = 'L7'; Value = 31 }
    @{ Sheet = 'Archer Heights'; Cell = 'L2'; Value = 16 }
    @{ Sheet = 'Archer Heights'; Cell = 'L7'; Value = 44 }
    @{ Sheet = 'Sauganash,Forest Glen'; Cell = 'L4'; Value = 1 }
    @{ Sheet = 'Sauganash,Forest Glen'; Cell = 'L7'; Value = 12 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
